# Fix table heading typo: "RMLSE" -> "RMSLE"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "RMSLE"

# Update selection to match the saved view (A2 selected instead of H6)
$ws.Range("A2").Select()
